$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 4708
$ws.Range("F3").Value = 1872
$ws.Range("F4").Value = 150
$ws.Range("F6").Value = 3167
$ws.Range("F8").Value = 596
$ws.Range("F10").Value = 648
$ws.Range("F11").Value = 552
$ws.Range("F12").Value = 545
$ws.Range("F13").Value = 407
$ws.Range("F15").Value = 1796
$ws.Range("F16").Value = 1383
$ws.Range("F18").Value = 1644
$ws.Range("F19").Value = 24
$ws.Range("F21").Value = 620
$ws.Range("F24").Value = 541
$ws.Range("F26").Value = 56
$ws.Range("F27").Value = 111
$ws.Range("F28").Value = 7
$ws.Range("F30").Value = 43
$ws.Range("F32").Value = 3992
$ws.Range("F33").Value = 10
$ws.Range("F35").Value = 85
$ws.Range("F36").Value = 1502
$ws.Range("F37").Value = 60
$ws.Range("F38").Value = 1891

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 28
$ws.Range("F3").Value = 57

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 4708
$ws.Range("F3").Value = 1872
$ws.Range("F4").Value = 150
$ws.Range("F6").Value = 3167
$ws.Range("F8").Value = 596
$ws.Range("F10").Value = 648
$ws.Range("F11").Value = 552
$ws.Range("F12").Value = 545
$ws.Range("F13").Value = 28
$ws.Range("F14").Value = 407
$ws.Range("F16").Value = 1796
$ws.Range("F17").Value = 1383
$ws.Range("F19").Value = 1644
$ws.Range("F20").Value = 24
$ws.Range("F22").Value = 620
$ws.Range("F25").Value = 541
$ws.Range("F27").Value = 56
$ws.Range("F28").Value = 111
$ws.Range("F29").Value = 7
$ws.Range("F31").Value = 43
$ws.Range("F33").Value = 3992
$ws.Range("F34").Value = 57
$ws.Range("F35").Value = 10
$ws.Range("F38").Value = 85
$ws.Range("F39").Value = 1502
$ws.Range("F40").Value = 60
$ws.Range("F41").Value = 1891
